$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 2 de Abril de 2020 a las 00:50"
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 212980
$ws.Cells.Item(4, 3).Value = 24450
$ws.Cells.Item(4, 4).Value = 8805
$ws.Cells.Item(4, 5).Value = 199416
$ws.Cells.Item(4, 6).Value = 5005
$ws.Cells.Item(4, 7).Value = 706
$ws.Cells.Item(4, 8).Value = 4759

$ws.Cells.Item(8, 1).Value = "Alemania"
$ws.Cells.Item(8, 2).Value = 77981
$ws.Cells.Item(8, 3).Value = 6173
$ws.Cells.Item(8, 4).Value = 18700
$ws.Cells.Item(8, 5).Value = 58350
$ws.Cells.Item(8, 6).Value = 3408
$ws.Cells.Item(8, 7).Value = 156
$ws.Cells.Item(8, 8).Value = 931

$ws.Cells.Item(25, 1).Value = "Chequia"
$ws.Cells.Item(25, 2).Value = 3589
$ws.Cells.Item(25, 3).Value = 281
$ws.Cells.Item(25, 4).Value = 61
$ws.Cells.Item(25, 5).Value = 3489
$ws.Cells.Item(25, 6).Value = 70
$ws.Cells.Item(25, 7).Value = 8
$ws.Cells.Item(25, 8).Value = 39

$ws.Cells.Item(31, 1).Value = "Ecuador"
$ws.Cells.Item(31, 2).Value = 2758
$ws.Cells.Item(31, 3).Value = 456
$ws.Cells.Item(31, 4).Value = 58
$ws.Cells.Item(31, 5).Value = 2602
$ws.Cells.Item(31, 6).Value = 100
$ws.Cells.Item(31, 7).Value = 19
$ws.Cells.Item(31, 8).Value = 98

$ws.Cells.Item(34, 1).Value = "Japon"
$ws.Cells.Item(34, 2).Value = 2384
$ws.Cells.Item(34, 3).Value = 206
$ws.Cells.Item(34, 4).Value = 472
$ws.Cells.Item(34, 5).Value = 1855
$ws.Cells.Item(34, 6).Value = 69
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 57

$ws.Cells.Item(35, 1).Value = "Luxemburgo"
$ws.Cells.Item(35, 2).Value = 2319
$ws.Cells.Item(35, 3).Value = 141
$ws.Cells.Item(35, 4).Value = 80
$ws.Cells.Item(35, 5).Value = 2210
$ws.Cells.Item(35, 6).Value = 31
$ws.Cells.Item(35, 7).Value = 6
$ws.Cells.Item(35, 8).Value = 29

$ws.Cells.Item(36, 1).Value = "Filipinas"
$ws.Cells.Item(36, 2).Value = 2311
$ws.Cells.Item(36, 3).Value = 227
$ws.Cells.Item(36, 4).Value = 50
$ws.Cells.Item(36, 5).Value = 2165
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 8
$ws.Cells.Item(36, 8).Value = 96

$ws.Cells.Item(52, 1).Value = "Argentina"
$ws.Cells.Item(52, 2).Value = 1054
$ws.Cells.Item(52, 3).Value = 0
$ws.Cells.Item(52, 4).Value = 248
$ws.Cells.Item(52, 5).Value = 775
$ws.Cells.Item(52, 6).Value = 0
$ws.Cells.Item(52, 7).Value = 4
$ws.Cells.Item(52, 8).Value = 31

$ws.Cells.Item(159, 1).Value = "Birmania"
$ws.Cells.Item(159, 2).Value = 16
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 4).Value = 0
$ws.Cells.Item(159, 5).Value = 15
$ws.Cells.Item(159, 6).Value = 0
$ws.Cells.Item(159, 7).Value = 0
$ws.Cells.Item(159, 8).Value = 1

$ws.Cells.Item(160, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(160, 2).Value = 16
$ws.Cells.Item(160, 3).Value = 10
$ws.Cells.Item(160, 4).Value = 6
$ws.Cells.Item(160, 5).Value = 9
$ws.Cells.Item(160, 6).Value = 0
$ws.Cells.Item(160, 7).Value = 1
$ws.Cells.Item(160, 8).Value = 1

$ws.Cells.Item(161, 1).Value = "Eritrea"
$ws.Cells.Item(161, 2).Value = 15
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 0
$ws.Cells.Item(161, 5).Value = 15
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 0

$ws.Cells.Item(172, 1).Value = "Mozambique"
$ws.Cells.Item(172, 2).Value = 10
$ws.Cells.Item(172, 3).Value = 2
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 10
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 0

$ws.Cells.Item(173, 1).Value = "Laos"
$ws.Cells.Item(173, 2).Value = 10
$ws.Cells.Item(173, 3).Value = 1
$ws.Cells.Item(173, 4).Value = 0
$ws.Cells.Item(173, 5).Value = 10
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 0

$ws.Cells.Item(174, 1).Value = "Surinam"
$ws.Cells.Item(174, 2).Value = 10
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(174, 4).Value = 0
$ws.Cells.Item(174, 5).Value = 10
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

$ws.Cells.Item(194, 1).Value = "Nicaragua"
$ws.Cells.Item(194, 2).Value = 5
$ws.Cells.Item(194, 3).Value = 0
$ws.Cells.Item(194, 4).Value = 0
$ws.Cells.Item(194, 5).Value = 4
$ws.Cells.Item(194, 6).Value = 0
$ws.Cells.Item(194, 7).Value = 0
$ws.Cells.Item(194, 8).Value = 1

$ws.Cells.Item(195, 1).Value = "Somalia"
$ws.Cells.Item(195, 2).Value = 5
$ws.Cells.Item(195, 3).Value = 0
$ws.Cells.Item(195, 4).Value = 1
$ws.Cells.Item(195, 5).Value = 4
$ws.Cells.Item(195, 6).Value = 0
$ws.Cells.Item(195, 7).Value = 0
$ws.Cells.Item(195, 8).Value = 0

$ws.Cells.Item(200, 1).Value = "Belice"
$ws.Cells.Item(200, 2).Value = 3
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 0
$ws.Cells.Item(200, 5).Value = 3
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 0

$ws.Cells.Item(201, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(201, 2).Value = 3
$ws.Cells.Item(201, 3).Value = 0
$ws.Cells.Item(201, 4).Value = 0
$ws.Cells.Item(201, 5).Value = 3
$ws.Cells.Item(201, 6).Value = 0
$ws.Cells.Item(201, 7).Value = 0
$ws.Cells.Item(201, 8).Value = 0

$ws.Cells.Item(202, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(202, 2).Value = 3
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 0
$ws.Cells.Item(202, 5).Value = 3
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

$ws.Cells.Item(204, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(204, 2).Value = 2
$ws.Cells.Item(204, 3).Value = 2
$ws.Cells.Item(204, 4).Value = 0
$ws.Cells.Item(204, 5).Value = 2
$ws.Cells.Item(204, 6).Value = 0
$ws.Cells.Item(204, 7).Value = 0
$ws.Cells.Item(204, 8).Value = 0

$ws.Cells.Item(205, 1).Value = "Anguila"
$ws.Cells.Item(205, 2).Value = 2
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 0
$ws.Cells.Item(205, 5).Value = 2
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

$ws.Cells.Item(206, 1).Value = "Burundi"
$ws.Cells.Item(206, 2).Value = 2
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 0
$ws.Cells.Item(206, 5).Value = 2
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

$ws.Cells.Item(207, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(207, 2).Value = 2
$ws.Cells.Item(207, 3).Value = 1
$ws.Cells.Item(207, 4).Value = 1
$ws.Cells.Item(207, 5).Value = 1
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

$ws.Cells.Item(208, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(208, 2).Value = 1
$ws.Cells.Item(208, 3).Value = 0
$ws.Cells.Item(208, 4).Value = 0
$ws.Cells.Item(208, 5).Value = 1
$ws.Cells.Item(208, 6).Value = 0
$ws.Cells.Item(208, 7).Value = 0
$ws.Cells.Item(208, 8).Value = 0

$ws.Cells.Item(209, 1).Value = "Timor Oriental"
$ws.Cells.Item(209, 2).Value = 1
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 0
$ws.Cells.Item(209, 5).Value = 1
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0
